$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Person")

$ws.Range("B1").Value = "foaf:familyName"
$ws.Range("C1").Value = "givenName"
